$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9 (ALC)
$ws.Range("H9").Value = 356.7143
$ws.Range("I9").Value = 299.5
$ws.Range("J9").Value = 399.625
$ws.Range("K9").Value = 299.5
$ws.Range("L9").Value = 399.625
$ws.Range("M9").Value = -130.5
$ws.Range("N9").Value = -737.625

# Row 19 (ALC)
$ws.Range("H19").Value = 16427.53
$ws.Range("I19").Value = 47669
$ws.Range("J19").Value = 3410.25
$ws.Range("K19").Value = 47669
$ws.Range("L19").Value = 3410.25
$ws.Range("M19").Value = -47494
$ws.Range("N19").Value = -3760.25

# Row 40 (ALC)
$ws.Range("H40").Value = 14333.6
$ws.Range("J40").Value = 12917.125
$ws.Range("L40").Value = 12917.125
$ws.Range("N40").Value = -13267.125

# Row 92 (ALC)
$ws.Range("H92").Value = 629.2105
$ws.Range("I92").Value = 637.3333
$ws.Range("K92").Value = 637.3333
$ws.Range("M92").Value = 610.6667

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 1090.1154
$ws.Range("I2").Value = 1081.8182
$ws.Range("J2").Value = 1135.75
$ws.Range("K2").Value = 1081.8182
$ws.Range("L2").Value = 1135.75
$ws.Range("M2").Value = -968.8181999999999
$ws.Range("N2").Value = -1361.75

# Row 116 (ARM)
$ws.Range("H116").Value = 1090.1154
$ws.Range("I116").Value = 1081.8182
$ws.Range("J116").Value = 1135.75
$ws.Range("K116").Value = 1081.8182
$ws.Range("L116").Value = 1135.75
$ws.Range("M116").Value = 1212.1818
$ws.Range("N116").Value = -5723.75

# Row 132 (ARM)
$ws.Range("H132").Value = 8422.406999999999
$ws.Range("I132").Value = 3172.55
$ws.Range("J132").Value = 23422
$ws.Range("K132").Value = 9517.650000000001
$ws.Range("L132").Value = 70266
$ws.Range("M132").Value = -6987.650000000001
$ws.Range("N132").Value = -75326

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 1090.1154
$ws.Range("I3").Value = 1081.8182
$ws.Range("J3").Value = 1135.75
$ws.Range("K3").Value = 1081.8182
$ws.Range("L3").Value = 1135.75
$ws.Range("M3").Value = -967.8181999999999
$ws.Range("N3").Value = -1363.75

# Row 94 (BSM)
$ws.Range("H94").Value = 2429.2307
$ws.Range("I94").Value = 1095.3235
$ws.Range("K94").Value = 1095.3235
$ws.Range("M94").Value = -644.3235

# Row 99 (BSM)
$ws.Range("H99").Value = 6429.65
$ws.Range("I99").Value = 6605.9375
$ws.Range("K99").Value = 6605.9375
$ws.Range("M99").Value = -5107.9375

$ws = $wb.Worksheets.Item("CRP")
# Row 105 (CRP)
$ws.Range("H105").Value = 1044.3077
$ws.Range("I105").Value = 1062.7727
$ws.Range("K105").Value = 1062.7727
$ws.Range("M105").Value = 684.2273

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws.Range("H5").Value = 588.3514
$ws.Range("I5").Value = 540.35297
$ws.Range("J5").Value = 1132.3334
$ws.Range("K5").Value = 1621.05891
$ws.Range("L5").Value = 3397.0002
$ws.Range("M5").Value = -1509.05891
$ws.Range("N5").Value = -3621.0002

# Row 23 (CUL)
$ws.Range("H23").Value = 902.7273
$ws.Range("I23").Value = 1443.25
$ws.Range("K23").Value = 4329.75
$ws.Range("M23").Value = -4094.75

# Row 74 (CUL)
$ws.Range("H74").Value = 3992.6667
$ws.Range("I74").Value = 3992.6667
$ws.Range("K74").Value = 11978.0001
$ws.Range("M74").Value = -10917.0001

# Row 77 (CUL)
$ws.Range("H77").Value = 3992.6667
$ws.Range("I77").Value = 3992.6667
$ws.Range("K77").Value = 35934.0003
$ws.Range("M77").Value = -30630.0003

# Row 122 (CUL)
$ws.Range("H122").Value = 788.4783
$ws.Range("I122").Value = 516
$ws.Range("J122").Value = 1085.7273
$ws.Range("K122").Value = 4644
$ws.Range("L122").Value = 9771.545700000001
$ws.Range("M122").Value = -2194
$ws.Range("N122").Value = -14671.5457

# Row 124 (CUL)
$ws.Range("H124").Value = 6105.5835
$ws.Range("I124").Value = 3721
$ws.Range("K124").Value = 11163
$ws.Range("M124").Value = -6253

# Row 135 (CUL)
$ws.Range("H135").Value = 588.3514
$ws.Range("I135").Value = 540.35297
$ws.Range("J135").Value = 1132.3334
$ws.Range("K135").Value = 4863.17673
$ws.Range("L135").Value = 10191.0006
$ws.Range("M135").Value = -2328.17673
$ws.Range("N135").Value = -15261.0006

$ws = $wb.Worksheets.Item("GSM")
# Row 97 (GSM)
$ws.Range("H97").Value = 648.6818
$ws.Range("I97").Value = 469.58823
$ws.Range("J97").Value = 1257.6
$ws.Range("K97").Value = 469.58823
$ws.Range("L97").Value = 1257.6
$ws.Range("M97").Value = 26.41176999999999
$ws.Range("N97").Value = -2249.6

# Row 113 (GSM)
$ws.Range("H113").Value = 6185.125
$ws.Range("I113").Value = 5995.75
$ws.Range("J113").Value = 6374.5
$ws.Range("K113").Value = 5995.75
$ws.Range("L113").Value = 6374.5
$ws.Range("M113").Value = -3825.75
$ws.Range("N113").Value = -10714.5

# Row 126 (GSM)
$ws.Range("H126").Value = 4197.275
$ws.Range("I126").Value = 4007.44
$ws.Range("J126").Value = 4513.6665
$ws.Range("K126").Value = 12022.32
$ws.Range("L126").Value = 13540.9995
$ws.Range("M126").Value = -9552.32
$ws.Range("N126").Value = -18480.9995

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (LTW)
$ws.Range("H16").Value = 1436.2894
$ws.Range("J16").Value = 1299.3334
$ws.Range("L16").Value = 1299.3334
$ws.Range("N16").Value = -1639.3334

# Row 55 (LTW)
$ws.Range("H55").Value = 644.4545000000001
$ws.Range("I55").Value = 241.83333
$ws.Range("K55").Value = 241.83333
$ws.Range("M55").Value = -68.83332999999999

# Row 61 (LTW)
$ws.Range("H61").Value = 28248
$ws.Range("I61").Value = 36997.668
$ws.Range("J61").Value = 1999
$ws.Range("K61").Value = 36997.668
$ws.Range("L61").Value = 1999
$ws.Range("M61").Value = -36795.668
$ws.Range("N61").Value = -2403

# Row 93 (LTW)
$ws.Range("H93").Value = 7065.2104
$ws.Range("I93").Value = 7118.8335
$ws.Range("K93").Value = 7118.8335
$ws.Range("M93").Value = -5870.8335

# Row 113 (LTW)
$ws.Range("H113").Value = 28248
$ws.Range("I113").Value = 36997.668
$ws.Range("J113").Value = 1999
$ws.Range("K113").Value = 36997.668
$ws.Range("L113").Value = 1999
$ws.Range("M113").Value = -34827.668
$ws.Range("N113").Value = -6339

$ws = $wb.Worksheets.Item("WVR")
# Row 100 (WVR)
$ws.Range("H100").Value = 977.5172
$ws.Range("I100").Value = 443.26315
$ws.Range("J100").Value = 1992.6
$ws.Range("K100").Value = 886.5263
$ws.Range("L100").Value = 3985.2
$ws.Range("M100").Value = -345.5263
$ws.Range("N100").Value = -5067.2

# Row 113 (WVR)
$ws.Range("H113").Value = 3730.818
$ws.Range("I113").Value = 3255
$ws.Range("K113").Value = 9765
$ws.Range("M113").Value = -7595

# Row 126 (WVR)
$ws.Range("H126").Value = 3128400.5
$ws.Range("I126").Value = 4169659.2
$ws.Range("K126").Value = 12508977.6
$ws.Range("M126").Value = -12508977.6
